$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.4375153333333333"
$ws.Range("H2").Value = [double]"1.312546"
$ws.Range("I2").Value = [double]"0.002535486401940996"
$ws.Range("J2").Value = [double]"0.002555908833496712"
$ws.Range("M2").Value = [double]"2.912114666666667"
$ws.Range("N2").Value = [double]"8.736344000000001"
$ws.Range("O2").Value = [double]"0.04564398277650125"
$ws.Range("P2").Value = [double]"0.06163513710720567"
$ws.Range("Q2").Value = [double]"1.274094819091556"
$ws.Range("R2").Value = [double]"11.466853371824"
$ws.Range("S2").Value = [double]"0.000115729697660248"
$ws.Range("T2").Value = [double]"0.0001575337913860879"
$ws.Range("G3").Value = [double]"0.4375153333333333"
$ws.Range("H3").Value = [double]"1.312546"
$ws.Range("I3").Value = [double]"0.002535486401940996"
$ws.Range("J3").Value = [double]"0.002555908833496712"
$ws.Range("O3").Value = [double]"0.01351577128599483"
$ws.Range("P3").Value = [double]"0.01825095808139687"
$ws.Range("Q3").Value = [double]"0.3772758888248889"
$ws.Range("R3").Value = [double]"3.395482999424"
$ws.Range("S3").Value = [double]"3.426905430738446E-05"
$ws.Range("T3").Value = [double]"4.664778498002045E-05"
$ws.Range("G4").Value = [double]"0.4375153333333333"
$ws.Range("H4").Value = [double]"1.312546"
$ws.Range("I4").Value = [double]"0.002535486401940996"
$ws.Range("J4").Value = [double]"0.002555908833496712"
$ws.Range("M4").Value = [double]"7.668087"
$ws.Range("N4").Value = [double]"23.004261"
$ws.Range("O4").Value = [double]"0.12018827244785"
$ws.Range("P4").Value = [double]"0.1622956674765719"
$ws.Range("Q4").Value = [double]"3.354905639834"
$ws.Range("R4").Value = [double]"30.194150758506"
$ws.Range("S4").Value = [double]"0.0003047357304643032"
$ws.Range("T4").Value = [double]"0.0004148129301416152"
$ws.Range("G5").Value = [double]"0.4375153333333333"
$ws.Range("H5").Value = [double]"1.312546"
$ws.Range("I5").Value = [double]"0.002535486401940996"
$ws.Range("J5").Value = [double]"0.002555908833496712"
$ws.Range("M5").Value = [double]"49.6589625"
$ws.Range("N5").Value = [double]"99.317925"
$ws.Range("O5").Value = [double]"0.7783460091712006"
$ws.Range("P5").Value = [double]"0.7006905777265834"
$ws.Range("Q5").Value = [double]"21.726557531175"
$ws.Range("R5").Value = [double]"130.35934518705"
$ws.Range("S5").Value = [double]"0.001973485722258621"
$ws.Range("T5").Value = [double]"0.001790901237159289"
$ws.Range("G6").Value = [double]"0.4375153333333333"
$ws.Range("H6").Value = [double]"1.312546"
$ws.Range("I6").Value = [double]"0.002535486401940996"
$ws.Range("J6").Value = [double]"0.002555908833496712"
$ws.Range("M6").Value = [double]"2.699147"
$ws.Range("N6").Value = [double]"8.097441"
$ws.Range("O6").Value = [double]"0.04230596431845346"
$ws.Range("P6").Value = [double]"0.05712765960824214"
$ws.Range("Q6").Value = [double]"1.180918199420667"
$ws.Range("R6").Value = [double]"10.628263794786"
$ws.Range("S6").Value = [double]"0.0001072661972504397"
$ws.Range("T6").Value = [double]"0.0001460130898296994"
$ws.Range("I7").Value = [double]"0.005820447907772805"
$ws.Range("J7").Value = [double]"0.005867329523437988"
$ws.Range("M7").Value = [double]"2.912114666666667"
$ws.Range("N7").Value = [double]"8.736344000000001"
$ws.Range("O7").Value = [double]"0.04564398277650125"
$ws.Range("P7").Value = [double]"0.06163513710720567"
$ws.Range("Q7").Value = [double]"2.924804691679112"
$ws.Range("R7").Value = [double]"26.323242225112"
$ws.Range("S7").Value = [double]"0.0002656684240539047"
$ws.Range("T7").Value = [double]"0.0003616336596302562"
$ws.Range("I8").Value = [double]"0.005820447907772805"
$ws.Range("J8").Value = [double]"0.005867329523437988"
$ws.Range("O8").Value = [double]"0.01351577128599483"
$ws.Range("P8").Value = [double]"0.01825095808139687"
$ws.Range("S8").Value = [double]"7.866784270350437E-05"
$ws.Range("T8").Value = [double]"0.000107084385182009"
$ws.Range("I9").Value = [double]"0.005820447907772805"
$ws.Range("J9").Value = [double]"0.005867329523437988"
$ws.Range("M9").Value = [double]"7.668087"
$ws.Range("N9").Value = [double]"23.004261"
$ws.Range("O9").Value = [double]"0.12018827244785"
$ws.Range("P9").Value = [double]"0.1622956674765719"
$ws.Range("Q9").Value = [double]"7.701501967117"
$ws.Range("R9").Value = [double]"69.313517704053"
$ws.Range("S9").Value = [double]"0.0006995495789079161"
$ws.Range("T9").Value = [double]"0.0009522421613113649"
$ws.Range("I10").Value = [double]"0.005820447907772805"
$ws.Range("J10").Value = [double]"0.005867329523437988"
$ws.Range("M10").Value = [double]"49.6589625"
$ws.Range("N10").Value = [double]"99.317925"
$ws.Range("O10").Value = [double]"0.7783460091712006"
$ws.Range("P10").Value = [double]"0.7006905777265834"
$ws.Range("Q10").Value = [double]"49.87535970558751"
$ws.Range("R10").Value = [double]"299.252158233525"
$ws.Range("S10").Value = [double]"0.004530322400603827"
$ws.Range("T10").Value = [double]"0.004111182513490004"
$ws.Range("I11").Value = [double]"0.005820447907772805"
$ws.Range("J11").Value = [double]"0.005867329523437988"
$ws.Range("M11").Value = [double]"2.699147"
$ws.Range("N11").Value = [double]"8.097441"
$ws.Range("O11").Value = [double]"0.04230596431845346"
$ws.Range("P11").Value = [double]"0.05712765960824214"
$ws.Range("Q11").Value = [double]"2.710908982910333"
$ws.Range("R11").Value = [double]"24.398180846193"
$ws.Range("S11").Value = [double]"0.0002462396615036534"
$ws.Range("T11").Value = [double]"0.000335186803824355"
$ws.Range("G12").Value = [double]"99.58055866666666"
$ws.Range("H12").Value = [double]"298.741676"
$ws.Range("I12").Value = [double]"0.577088694179909"
$ws.Range("J12").Value = [double]"0.5817369361698658"
$ws.Range("M12").Value = [double]"2.912114666666667"
$ws.Range("N12").Value = [double]"8.736344000000001"
$ws.Range("O12").Value = [double]"0.04564398277650125"
$ws.Range("P12").Value = [double]"0.06163513710720567"
$ws.Range("Q12").Value = [double]"289.9900054080605"
$ws.Range("R12").Value = [double]"2609.910048672544"
$ws.Range("S12").Value = [double]"0.02634062641766137"
$ws.Range("T12").Value = [double]"0.03585543582115543"
$ws.Range("G13").Value = [double]"99.58055866666666"
$ws.Range("H13").Value = [double]"298.741676"
$ws.Range("I13").Value = [double]"0.577088694179909"
$ws.Range("J13").Value = [double]"0.5817369361698658"
$ws.Range("O13").Value = [double]"0.01351577128599483"
$ws.Range("P13").Value = [double]"0.01825095808139687"
$ws.Range("Q13").Value = [double]"85.86977625312711"
$ws.Range("R13").Value = [double]"772.827986278144"
$ws.Range("S13").Value = [double]"0.007799798802269067"
$ws.Range("T13").Value = [double]"0.01061725643643647"
$ws.Range("G14").Value = [double]"99.58055866666666"
$ws.Range("H14").Value = [double]"298.741676"
$ws.Range("I14").Value = [double]"0.577088694179909"
$ws.Range("J14").Value = [double]"0.5817369361698658"
$ws.Range("M14").Value = [double]"7.668087"
$ws.Range("N14").Value = [double]"23.004261"
$ws.Range("O14").Value = [double]"0.12018827244785"
$ws.Range("P14").Value = [double]"0.1622956674765719"
$ws.Range("Q14").Value = [double]"763.592387364604"
$ws.Range("R14").Value = [double]"6872.331486281436"
$ws.Range("S14").Value = [double]"0.06935929320266887"
$ws.Range("T14").Value = [double]"0.09441338435146428"
$ws.Range("G15").Value = [double]"99.58055866666666"
$ws.Range("H15").Value = [double]"298.741676"
$ws.Range("I15").Value = [double]"0.577088694179909"
$ws.Range("J15").Value = [double]"0.5817369361698658"
$ws.Range("M15").Value = [double]"49.6589625"
$ws.Range("N15").Value = [double]"99.317925"
$ws.Range("O15").Value = [double]"0.7783460091712006"
$ws.Range("P15").Value = [double]"0.7006905777265834"
$ws.Range("Q15").Value = [double]"4945.06722855705"
$ws.Range("R15").Value = [double]"29670.4033713423"
$ws.Range("S15").Value = [double]"0.4491746820527516"
$ws.Range("T15").Value = [double]"0.4076175898897558"
$ws.Range("G16").Value = [double]"99.58055866666666"
$ws.Range("H16").Value = [double]"298.741676"
$ws.Range("I16").Value = [double]"0.577088694179909"
$ws.Range("J16").Value = [double]"0.5817369361698658"
$ws.Range("M16").Value = [double]"2.699147"
$ws.Range("N16").Value = [double]"8.097441"
$ws.Range("O16").Value = [double]"0.04230596431845346"
$ws.Range("P16").Value = [double]"0.05712765960824214"
$ws.Range("Q16").Value = [double]"268.7825661834573"
$ws.Range("R16").Value = [double]"2419.043095651116"
$ws.Range("S16").Value = [double]"0.02441429370455813"
$ws.Range("T16").Value = [double]"0.03323326967105378"
$ws.Range("G17").Value = [double]"4.1363315"
$ws.Range("H17").Value = [double]"8.272663"
$ws.Range("I17").Value = [double]"0.02397084507248554"
$ws.Range("J17").Value = [double]"0.01610928107528529"
$ws.Range("M17").Value = [double]"2.912114666666667"
$ws.Range("N17").Value = [double]"8.736344000000001"
$ws.Range("O17").Value = [double]"0.04564398277650125"
$ws.Range("P17").Value = [double]"0.06163513710720567"
$ws.Range("Q17").Value = [double]"12.04547162734533"
$ws.Range("R17").Value = [double]"72.272829764072"
$ws.Range("S17").Value = [double]"0.00109412483962671"
$ws.Range("T17").Value = [double]"0.0009928977477737227"
$ws.Range("G18").Value = [double]"4.1363315"
$ws.Range("H18").Value = [double]"8.272663"
$ws.Range("I18").Value = [double]"0.02397084507248554"
$ws.Range("J18").Value = [double]"0.01610928107528529"
$ws.Range("O18").Value = [double]"0.01351577128599483"
$ws.Range("P18").Value = [double]"0.01825095808139687"
$ws.Range("Q18").Value = [double]"3.566819318645333"
$ws.Range("R18").Value = [double]"21.400915911872"
$ws.Range("S18").Value = [double]"0.0003239844595317308"
$ws.Range("T18").Value = [double]"0.0002940098136264717"
$ws.Range("G19").Value = [double]"4.1363315"
$ws.Range("H19").Value = [double]"8.272663"
$ws.Range("I19").Value = [double]"0.02397084507248554"
$ws.Range("J19").Value = [double]"0.01610928107528529"
$ws.Range("M19").Value = [double]"7.668087"
$ws.Range("N19").Value = [double]"23.004261"
$ws.Range("O19").Value = [double]"0.12018827244785"
$ws.Range("P19").Value = [double]"0.1622956674765719"
$ws.Range("Q19").Value = [double]"31.7177498028405"
$ws.Range("R19").Value = [double]"190.306498817043"
$ws.Range("S19").Value = [double]"0.002881014458377094"
$ws.Range("T19").Value = [double]"0.002614466524681135"
$ws.Range("G20").Value = [double]"4.1363315"
$ws.Range("H20").Value = [double]"8.272663"
$ws.Range("I20").Value = [double]"0.02397084507248554"
$ws.Range("J20").Value = [double]"0.01610928107528529"
$ws.Range("M20").Value = [double]"49.6589625"
$ws.Range("N20").Value = [double]"99.317925"
$ws.Range("O20").Value = [double]"0.7783460091712006"
$ws.Range("P20").Value = [double]"0.7006905777265834"
$ws.Range("Q20").Value = [double]"205.4059308460687"
$ws.Range("R20").Value = [double]"821.623723384275"
$ws.Range("S20").Value = [double]"0.01865761159863026"
$ws.Range("T20").Value = [double]"0.01128762146340157"
$ws.Range("G21").Value = [double]"4.1363315"
$ws.Range("H21").Value = [double]"8.272663"
$ws.Range("I21").Value = [double]"0.02397084507248554"
$ws.Range("J21").Value = [double]"0.01610928107528529"
$ws.Range("M21").Value = [double]"2.699147"
$ws.Range("N21").Value = [double]"8.097441"
$ws.Range("O21").Value = [double]"0.04230596431845346"
$ws.Range("P21").Value = [double]"0.05712765960824214"
$ws.Range("Q21").Value = [double]"11.1645667592305"
$ws.Range("R21").Value = [double]"66.98740055538299"
$ws.Range("S21").Value = [double]"0.001014109716319749"
$ws.Range("T21").Value = [double]"0.0009202855258023951"
$ws.Range("G22").Value = [double]"67.39800266666667"
$ws.Range("H22").Value = [double]"202.194008"
$ws.Range("I22").Value = [double]"0.3905845264378918"
$ws.Range("J22").Value = [double]"0.3937305443979143"
$ws.Range("M22").Value = [double]"2.912114666666667"
$ws.Range("N22").Value = [double]"8.736344000000001"
$ws.Range("O22").Value = [double]"0.04564398277650125"
$ws.Range("P22").Value = [double]"0.06163513710720567"
$ws.Range("Q22").Value = [double]"196.2707120696392"
$ws.Range("R22").Value = [double]"1766.436408626752"
$ws.Range("S22").Value = [double]"0.01782783339749903"
$ws.Range("T22").Value = [double]"0.02426763608726018"
$ws.Range("G23").Value = [double]"67.39800266666667"
$ws.Range("H23").Value = [double]"202.194008"
$ws.Range("I23").Value = [double]"0.3905845264378918"
$ws.Range("J23").Value = [double]"0.3937305443979143"
$ws.Range("O23").Value = [double]"0.01351577128599483"
$ws.Range("P23").Value = [double]"0.01825095808139687"
$ws.Range("Q23").Value = [double]"58.11828620350578"
$ws.Range("R23").Value = [double]"523.064575831552"
$ws.Range("S23").Value = [double]"0.005279051127183147"
$ws.Range("T23").Value = [double]"0.007185959661171902"
$ws.Range("G24").Value = [double]"67.39800266666667"
$ws.Range("H24").Value = [double]"202.194008"
$ws.Range("I24").Value = [double]"0.3905845264378918"
$ws.Range("J24").Value = [double]"0.3937305443979143"
$ws.Range("M24").Value = [double]"7.668087"
$ws.Range("N24").Value = [double]"23.004261"
$ws.Range("O24").Value = [double]"0.12018827244785"
$ws.Range("P24").Value = [double]"0.1622956674765719"
$ws.Range("Q24").Value = [double]"516.813748074232"
$ws.Range("R24").Value = [double]"4651.323732668088"
$ws.Range("S24").Value = [double]"0.04694367947743179"
$ws.Range("T24").Value = [double]"0.06390076150897353"
$ws.Range("G25").Value = [double]"67.39800266666667"
$ws.Range("H25").Value = [double]"202.194008"
$ws.Range("I25").Value = [double]"0.3905845264378918"
$ws.Range("J25").Value = [double]"0.3937305443979143"
$ws.Range("M25").Value = [double]"49.6589625"
$ws.Range("N25").Value = [double]"99.317925"
$ws.Range("O25").Value = [double]"0.7783460091712006"
$ws.Range("P25").Value = [double]"0.7006905777265834"
$ws.Range("Q25").Value = [double]"3346.9148869989"
$ws.Range("R25").Value = [double]"20081.4893219934"
$ws.Range("S25").Value = [double]"0.3040099073969563"
$ws.Range("T25").Value = [double]"0.2758832826227767"
$ws.Range("G26").Value = [double]"67.39800266666667"
$ws.Range("H26").Value = [double]"202.194008"
$ws.Range("I26").Value = [double]"0.3905845264378918"
$ws.Range("J26").Value = [double]"0.3937305443979143"
$ws.Range("M26").Value = [double]"2.699147"
$ws.Range("N26").Value = [double]"8.097441"
$ws.Range("O26").Value = [double]"0.04230596431845346"
$ws.Range("P26").Value = [double]"0.05712765960824214"
$ws.Range("Q26").Value = [double]"181.9171167037253"
$ws.Range("R26").Value = [double]"1637.254050333528"
$ws.Range("S26").Value = [double]"0.01652405503882149"
$ws.Range("T26").Value = [double]"0.02249290451773192"
